$d = $word.ActiveDocument

# The conclusion paragraph originally ends with "... but that doesn't
# affect the ability of the program." The edit changes "ability" to
# "efficiency" so the sentence reads "... affect the efficiency of the
# program."
$d.Content.Find.Execute("ability of the program", $true, $false, $false, $false, $false,
                         $true, 1, $false, "efficiency of the program", 2)
